# Add a "description" column to the source table, inserted between "name" and "path".
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new blank column before the current "path" column (D); this shifts
# "path" to E and "last_update" to F, carrying their existing formatting/width
# with them.
$ws.Range("D1").EntireColumn.Insert()

# Header + data for the new "description" column. The description mirrors the
# "name" column for every existing row, except the "pdf_wiki" row which gets
# a distinct, more readable value.
$ws.Range("D1").Value = "description"
$ws.Range("D2").Value = "pdf wiki"
$ws.Range("D3").Value = $ws.Range("C3").Value2
$ws.Range("D4").Value = $ws.Range("C4").Value2
$ws.Range("D5").Value = $ws.Range("C5").Value2
$ws.Range("D6").Value = $ws.Range("C6").Value2
$ws.Range("D7").Value = $ws.Range("C7").Value2
$ws.Range("D8").Value = $ws.Range("C8").Value2

# New column width - matches the "name" column's width.
$ws.Columns.Item(4).ColumnWidth = $ws.Columns.Item(3).ColumnWidth

# Rebuild the table over the new range so its column headers/names line up
# with the shifted data (a plain Resize keeps stale cached column names).
$oldTable = $ws.ListObjects.Item("Tableau1")
$oldTable.Unlist()
$table = $ws.ListObjects.Add(1, $ws.Range("A1:F8"), 0, 1)
$table.Name = "Tableau1"
$table.TableStyle = "TableStyleMedium9"

# Leave the selection where the edit finished.
$ws.Range("D3").Select()
